$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "MassSpectrometryResults" sheet right before "Array".
#    It shares the same DataEntity-style column layout used by
#    DataEntity / AlignmentSet / VariantSet / Array.
# ---------------------------------------------------------------------
$arraySheet = $wb.Worksheets.Item("Array")
$msSheet = $wb.Worksheets.Add($arraySheet)
$msSheet.Name = "MassSpectrometryResults"
$msSheet.Outline.SummaryRow = 1
$msSheet.Outline.SummaryColumn = 1

$msSheet.Range("A1").Value = "data_path"
$msSheet.Range("B1").Value = "data_format"
$msSheet.Range("C1").Value = "has_sample"
$msSheet.Range("D1").Value = "has_reference"
$msSheet.Range("E1").Value = "id"
$msSheet.Range("F1").Value = "name"
$msSheet.Range("G1").Value = "description"

$msValidation = $msSheet.Range("B2:B1048576").Validation
$msValidation.Add(3, 1, 1, '"CRAM,FASTQ,Zarr,FASTA,VCF,BCF,mzTab"')
$msValidation.IgnoreBlank = $true
$msValidation.InCellDropdown = $true
$msValidation.ShowInput = $false
$msValidation.ShowError = $false

# Match the page-margin defaults used throughout the rest of this workbook
# (0.75in sides, 1in top/bottom, 0.5in header/footer).
$msSheet.PageSetup.LeftMargin = 54
$msSheet.PageSetup.RightMargin = 54
$msSheet.PageSetup.TopMargin = 72
$msSheet.PageSetup.BottomMargin = 72
$msSheet.PageSetup.HeaderMargin = 36
$msSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Add the new "mzTab" option to the data_format dropdown on every
#    other sheet that already offers that same CRAM/FASTQ/... list.
# ---------------------------------------------------------------------
$sheetsWithFormatList = @("DataEntity", "AlignmentSet", "VariantSet", "Array")
foreach ($sheetName in $sheetsWithFormatList) {
    $ws = $wb.Worksheets.Item($sheetName)
    $validation = $ws.Range("B2:B1048576").Validation
    $validation.Modify(3, 1, 1, '"CRAM,FASTQ,Zarr,FASTA,VCF,BCF,mzTab"')
}

# ---------------------------------------------------------------------
# 3. Assay sheet: insert a new "sample_processing" column between
#    "omics_type" (C) and "id" (formerly D, now shifts to E).
# ---------------------------------------------------------------------
$assay = $wb.Worksheets.Item("Assay")
$assay.Columns.Item(4).Insert()
$assay.Range("D1").Value = "sample_processing"

# ---------------------------------------------------------------------
# Restore the originally active sheet (inserting a worksheet makes it
# active, same as interactive Excel would do).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("NamedThing").Activate()
